$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the duplicate screen assignment: Austin Wadlow's screen number
# duplicated Richard Li's (both "12"). Reassign Austin Wadlow to screen 10.
$ws.Range("C3").Value = 10

# Move active selection to C4 (matches post-edit selection state)
$ws.Range("C4").Select()
